# SpaceSolo.docx edit
#
# Target change (from the OOXML diff):
#   - The first paragraph of the document (currently empty, centered,
#     sz/szCs = 56) gains a run with the text "mkk" using that same
#     sz/szCs = 56 formatting.
#   - The "_GoBack" bookmark (a single bookmarkStart/bookmarkEnd pair
#     that Word drops at the location of the user's last edit) moves
#     from the end of the document (right after the final ".") to sit
#     immediately after the new "mkk" run in the first paragraph.

$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1)
$r1 = $p1.Range

# Word's Font.Size setter only ever emits <w:sz>, never the matching
# <w:szCs> the diff requires, so plain "type text then set Font.Size"
# can't reproduce the target run properties. Instead we borrow the
# formatting from the "SpaceSolo" run in the 3rd paragraph (which
# already carries sz=56/szCs=56) via Range.FormattedText, then edit
# the copied text down to "mkk" in place -- replacing text inside an
# existing run (rather than inserting fresh text) preserves its rPr,
# including szCs.

# 1) Type a 9-character placeholder at the very start of paragraph 1
#    (same length as "SpaceSolo") so we have a same-sized destination
#    range to receive the borrowed formatting.
$r1.InsertAfter("123456789")

# 2) Grab the fully-formatted "SpaceSolo" range (paragraph 3) and copy
#    its formatted text onto the placeholder range in paragraph 1.
$p3 = $d.Paragraphs(3)
$srcStart = $p3.Range.Start
$srcEnd = $srcStart + 9
$srcRange = $d.Range($srcStart, $srcEnd)

$destRange = $d.Range(0, 9)
$destRange.FormattedText = $srcRange.FormattedText

# 3) In-place replace the first 3 characters with "mkk" -- this is a
#    same-run text edit, so the run keeps its sz/szCs formatting.
$sub = $d.Range(0, 3)
$sub.Text = "mkk"

# 4) Remove the old "_GoBack" bookmark (currently sitting at the end
#    of the document, right after the final "." run) before we add
#    the new one, since a document can only have one bookmark of a
#    given name at a time.
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBookmark = $d.Bookmarks.Item("_GoBack")
    $oldBookmark.Delete()
}

# 5) Add the "_GoBack" bookmark right after "mkk" (position 3), i.e.
#    before the leftover placeholder characters that still follow it.
$bookmarkRange = $d.Range(3, 3)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# 6) Finally, delete the unused remainder of the placeholder text
#    ("ceSolo"), leaving the run's visible text as just "mkk".
$remainder = $d.Range(3, 9)
$remainder.Delete()

Write-Output ("Paragraph 1 text: [" + $d.Paragraphs(1).Range.Text + "]")
$newBookmark = $d.Bookmarks.Item("_GoBack")
Write-Output ("_GoBack bookmark at: " + $newBookmark.Start + "-" + $newBookmark.End)
